$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text cells stay text (avoid Excel auto-converting numeric-looking strings to numbers)
$ws.Range("B2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "30.151.89"
$ws.Range("E2").Value = "  +0.96%  "

$ws.Range("D3").Value = "1.895.92"
$ws.Range("E3").Value = "  +0.49%  "

$ws.Range("D4").Value = "0.9967"
$ws.Range("E4").Value = "  -0.48%  "

$ws.Range("D5").Value = "0.7452"
$ws.Range("E5").Value = "  -0.41%  "

$ws.Range("D6").Value = "243.32"
$ws.Range("E6").Value = "  +0.36%  "

$ws.Range("D7").Value = "0.9991"
$ws.Range("E7").Value = "  -0.23%  "

$ws.Range("D8").Value = "0.3172"
$ws.Range("E8").Value = "  +1.77%  "

$ws.Range("D9").Value = "0.07253"
$ws.Range("E9").Value = "  +1.79%  "

$ws.Range("D10").Value = "25.10"
$ws.Range("E10").Value = "  -0.96%  "

$ws.Range("D11").Value = "0.08365"
$ws.Range("E11").Value = "  -1.87%  "

$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.949.49"
$ws.Range("E12").Value = "  +0.81%  "

$ws.Range("B13").Value = "Polygon"
$ws.Range("C13").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D13").Value = "0.7630"
$ws.Range("E13").Value = "  +0.45%  "

$ws.Range("D14").Value = "5.436"
$ws.Range("E14").Value = "  +1.47%  "

$ws.Range("D15").Value = "93.23"
$ws.Range("E15").Value = "  -0.21%  "

$ws.Range("D16").Value = "6.176"
$ws.Range("E16").Value = "  +0.55%  "

$ws.Range("D17").Value = "30.161.20"
$ws.Range("E17").Value = "  +0.89%  "

$ws.Range("D18").Value = "250.14"
$ws.Range("E18").Value = "  +2.94%  "

$ws.Range("D19").Value = "13.69"
$ws.Range("E19").Value = "  -0.27%  "

$ws.Range("D20").Value = "0.000007882"
$ws.Range("E20").Value = "  +1.20%  "

$ws.Range("D21").Value = "2.166.40"
$ws.Range("E21").Value = "  +0.55%  "

$ws.Range("D22").Value = "0.9988"
$ws.Range("E22").Value = "  -0.15%  "

$ws.Range("D23").Value = "8.037"
$ws.Range("E23").Value = "  +0.70%  "

$ws.Range("D24").Value = "0.9982"
$ws.Range("E24").Value = "  -0.37%  "

$ws.Range("D25").Value = "0.1584"
$ws.Range("E25").Value = "  -0.74%  "

$ws.Range("D26").Value = "9.322"
$ws.Range("E26").Value = "  -0.37%  "

$ws.Range("D27").Value = "164.30"
$ws.Range("E27").Value = "  +0.98%  "

$ws.Range("D28").Value = "18.80"
$ws.Range("E28").Value = "  +0.31%  "

$ws.Range("D29").Value = "2.058"
$ws.Range("E29").Value = "  +1.57%  "

$ws.Range("D30").Value = "1.477"
$ws.Range("E30").Value = "  -2.49%  "

$ws.Range("D31").Value = "4.613"
$ws.Range("E31").Value = "  +3.07%  "

$ws.Range("D32").Value = "1.540"
$ws.Range("E32").Value = "  +0.65%  "

$ws.Range("D33").Value = "4.249"
$ws.Range("E33").Value = "  +3.74%  "

$ws.Range("D34").Value = "0.05397"
$ws.Range("E34").Value = "  +0.07%  "

$ws.Range("D35").Value = "1.258"
$ws.Range("E35").Value = "  +1.82%  "

$ws.Range("D36").Value = "0.7643"
$ws.Range("E36").Value = "  +2.97%  "

$ws.Range("D37").Value = "1.002"
$ws.Range("E37").Value = "  -0.12%  "

$ws.Range("D38").Value = "2.720"
$ws.Range("E38").Value = "  +0.32%  "

$ws.Range("D39").Value = "0.01974"
$ws.Range("E39").Value = "  +1.87%  "

$ws.Range("D40").Value = "2.768"
$ws.Range("E40").Value = "  -0.13%  "

$ws.Range("D41").Value = "0.4583"
$ws.Range("E41").Value = "  +3.00%  "

$ws.Range("D42").Value = "1.104.80"
$ws.Range("E42").Value = "  +0.36%  "

$ws.Range("D43").Value = "6.090"
$ws.Range("E43").Value = "  +0.31%  "

$ws.Range("D44").Value = "73.00"
$ws.Range("E44").Value = "  +0.86%  "

$ws.Range("D45").Value = "0.8719"
$ws.Range("E45").Value = "  +1.87%  "

$ws.Range("E46").Value = "  +2.16%  "

$ws.Range("D47").Value = "0.9997"
$ws.Range("E47").Value = "  -0.15%  "

$ws.Range("D48").Value = "1.876"
$ws.Range("E48").Value = "  +0.79%  "

$ws.Range("D49").Value = "7.643"
$ws.Range("E49").Value = "  -0.18%  "

$ws.Range("D50").Value = "9.601"
$ws.Range("E50").Value = "  -1.11%  "

$ws.Range("D51").Value = "2.053.06"
$ws.Range("E51").Value = "  +0.42%  "
